# The deck currently carries two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (not applied to any slide; only
#                            wired to the Notes Master relationship)
#   ppt/theme/theme2.xml -> "Integral"     (the theme actually applied to the
#                            slide master / whole presentation design)
#
# The target edit swaps which theme's color values live in the part that is
# actually applied to the presentation: after the edit, the live design
# (reached through SlideMaster / Slides / Design) must carry the "Office
# Theme" palette instead of "Integral". (The Notes-Master-only theme part
# is not reachable through the Slide/Design/ColorScheme object model, so it
# cannot be touched from here — only the applied design's 12 theme colors
# are settable through PowerPoint's COM surface.)
#
# PowerPoint doesn't expose "load this raw theme XML" on the object model;
# the supported, scriptable way to repaint a design's theme colors is to
# walk the 12-slot ThemeColorScheme on a slide (order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) and assign each RGB value in turn.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Target palette == the "Office Theme" colors that currently sit in
# ppt/theme/theme1.xml, in clrScheme slot order.
$officeThemeRgbHex = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeRgbHex[$i - 1]
    $red   = ($hex -shr 16) -band 0xFF
    $green = ($hex -shr 8)  -band 0xFF
    $blue  =  $hex          -band 0xFF

    # OLE RGB() packs as 0x00BBGGRR (little-endian B/G/R).
    $themeColors.Colors($i).RGB = ($blue * 65536) + ($green * 256) + $red
}
